$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111486450
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2"
$ws.Range("J2").ClearContents()
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "fullt utvecklade blad"
$ws.Range("Q2").Value = 624051.1502826829
$ws.Range("R2").Value = 6932945.755648845
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "13:43"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "13:43"
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("A3").Value = 111486415
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "10"
$ws.Range("Q3").Value = 624040.2038791699
$ws.Range("R3").Value = 6932953.67081845
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-08-12"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "13:46"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-08-12"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "13:46"

# Row 4
$ws.Range("A4").Value = 111486385
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "20"
$ws.Range("Q4").Value = 624029.7289273632
$ws.Range("R4").Value = 6932998.597210908
$ws.Range("S4").Value = 10
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "14:11"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "14:11"
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = "Ca 20 ex varav 3 blommande"

# Row 5
$ws.Range("A5").Value = 111486117
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "10"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "blomning"
$ws.Range("Q5").Value = 623993.6707231236
$ws.Range("R5").Value = 6933021.760048959
$ws.Range("S5").Value = 15
$ws.Range("AC5").NumberFormat = "@"
$ws.Range("AC5").Value = "10 plantor varav 2 blommande"

# Row 6
$ws.Range("A6").Value = 111485854
$ws.Range("B6").Value = 89405
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1202
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = ""
$ws.Range("N6").ClearContents()
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "Glödenhöjden nordost (Glödenhöjden), Mpd"
$ws.Range("Q6").Value = 624096.1730324102
$ws.Range("R6").Value = 6933042.231978768
$ws.Range("S6").Value = 20
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "15:02"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "15:02"
$ws.Range("AH6").NumberFormat = "@"
$ws.Range("AH6").Value = "Häll- och blockskog"

# Row 8
$ws.Range("A8").Value = 111486400
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "5"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "plantor/tuvor"
$ws.Range("Q8").Value = 624030.1824148977
$ws.Range("R8").Value = 6932961.620511409
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-08-14"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-08-14"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "00:00"

# Row 9
$ws.Range("A9").Value = 111486280
$ws.Range("B9").Value = 96348
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "3"
$ws.Range("J9").ClearContents()
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "blomning"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "observerad"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "Glödenhöjden (Glödenhöjden), Mpd"
$ws.Range("Q9").Value = 624009.7035872869
$ws.Range("R9").Value = 6933014.034667149
$ws.Range("S9").Value = 10
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "14:17"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "14:17"
$ws.Range("AH9").ClearContents()

# Row 10
$ws.Range("A10").Value = 111486347
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "20"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "plantor/tuvor"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "i frukt"
$ws.Range("Q10").Value = 624030.5648888731
$ws.Range("R10").Value = 6933013.425735661
$ws.Range("AC10").NumberFormat = "@"
$ws.Range("AC10").Value = "Ca 20 ex, flesta som bladrosetter. 1 överblommad fruktbildande"
